# Applies:
#  - Update BOM date from "2024 October 14" to "2024 October 29" on all sheets
#  - Fix connector part numbers on "All" and "Top" sheets (CT3151V1-x -> CT3149-x)

$wb = $excel.ActiveWorkbook

$oldDate = "2024 October 14"
$newDate = "2024 October 29"

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A2").Value() -eq $oldDate) {
        $ws.Range("A2").Value = $newDate
    }
}

$partNumberMap = @{
    "CT3151V1-0" = "CT3149-0"
    "CT3151V1-2" = "CT3149-2"
    "CT3151V1-4" = "CT3149-4"
}

foreach ($sheetName in @("All", "Top")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowNum in 5..7) {
        $cell = $ws.Cells.Item($rowNum, 4)
        $current = $cell.Value()
        if ($partNumberMap.ContainsKey($current)) {
            $cell.Value = $partNumberMap[$current]
        }
    }
}
